# TimeTrack.xlsx edit: add Monday 9.9.18 entry to row 3, update view zoom/selection.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Monday 9.9.18"
$ws.Range("B3").Value = "0900-1100"
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = "Setting up repo on laptop + building simple classes"

$ws.Range("B3").Select()
$excel.ActiveWindow.Zoom = 85
